# dev board: fix wrong parts on BOM
# J4's DigiKey part # was wrong -> correct it to " A122210-ND"
# J5's DigiKey part # was wrong -> correct it to " 609-3486-2-ND"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E11").Value = " A122210-ND"
$ws.Range("E12").Value = " 609-3486-2-ND"

# leave the cursor where the editor ended up after fixing the two rows
$ws.Range("E13").Select() | Out-Null
